$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-01 Sunday" "2023-10-02 Monday"

Replace-Text "31×70=" "48×35="
Replace-Text "53×89=" "40×16="
Replace-Text "26×60=" "90×15="
Replace-Text "24×89=" "81×24="
Replace-Text "33×95=" "21×61="

Replace-Text "98×12=" "48×50="
Replace-Text "61×30=" "60×83="
Replace-Text "75×87=" "39×22="
Replace-Text "18×91=" "87×67="
Replace-Text "74×38=" "15×82="

Replace-Text "53×93=" "90×71="
Replace-Text "35×99=" "43×87="
Replace-Text "72×64=" "15×11="
Replace-Text "48×11=" "12×51="
Replace-Text "34×11=" "18×24="

Replace-Text "59×88=" "63×89="
Replace-Text "19×40=" "91×34="
Replace-Text "19×87=" "11×93="
Replace-Text "29×57=" "43×24="
Replace-Text "59×70=" "30×49="

Replace-Text "28×97=" "64×59="
Replace-Text "97×23=" "21×48="
Replace-Text "72×37=" "97×27="
Replace-Text "87×23=" "60×94="
Replace-Text "26×14=" "61×11="
